$wb = $excel.ActiveWorkbook

# Both "展览" (sheet 1) and "全部类型" (sheet 4) contain the same set of
# listings (row numbers line up for most rows, except sheet 4 has two
# extra rows inserted before the tail, shifting rows 39/40 -> 41/43).

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # F5: 想去人数 16015 -> 16023
    $ws.Range("F5").Value = 16023

    # F9: 想去人数 15521 -> 15524
    $ws.Range("F9").Value = 15524

    # F11: 想去人数 9140 -> 9145
    $ws.Range("F11").Value = 9145

    # F18: 想去人数 213 -> 214
    $ws.Range("F18").Value = 214

    # D19: venue text updated
    $ws.Range("D19").Value = "娄东街道常胜北路66号1幢1楼 鑫锐体育迎篮而上篮球运动中心"

    # I19: cover image updated
    $ws.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202409/08LBKX1Q1727581256319.jpeg"

    # I21: cover image updated
    $ws.Range("I21").Value = "//i1.hdslb.com/bfs/openplatform/202409/w6kKyFT11727579615000.jpeg"

    # C22: event renamed (cancelled)
    $ws.Range("C22").Value = "昆山·心动次元动漫游戏嘉年华（取消）"

    # G22: 最低票价 45 -> 不可售 (now text, not sellable)
    $ws.Range("G22").Value = "不可售"

    # F25: 想去人数 1127 -> 1128
    $ws.Range("F25").Value = 1128

    # F26: 想去人数 10 -> 11
    $ws.Range("F26").Value = 11
}

# "展览" sheet specific rows (37, 39, 40)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F37").Value = 466
$wsExhibit.Range("F39").Value = 5625
$wsExhibit.Range("F40").Value = 5235
$wsExhibit.Range("G40").Value = 55
$wsExhibit.Range("I40").Value = "//i0.hdslb.com/bfs/openplatform/202409/6zugizcG1727576290688.jpeg"

# "全部类型" sheet specific rows (39, 41, 43) - shifted by the two extra
# rows present only on this sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F39").Value = 466
$wsAll.Range("F41").Value = 5625
$wsAll.Range("F43").Value = 5235
$wsAll.Range("G43").Value = 55
$wsAll.Range("I43").Value = "//i0.hdslb.com/bfs/openplatform/202409/6zugizcG1727576290688.jpeg"
